$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = $false
$ws.Range("E2").Value = 19.170000000000002
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = $false
$ws.Range("C3").Value = 0
